# Update "Agile_Template_Filled.xlsx" with the latest Agile Document data.
# The Product Backlog's Assignee column is consolidated onto "Team A"
# for the User Story rows that used to read "Team B" / "Both Teams".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# US002 (row 3) and US004 (row 5) were assigned to "Team B"; US006 (row 7)
# was assigned to "Both Teams". All three now belong to "Team A".
$ws.Range("G3").Value = "Team A"
$ws.Range("G5").Value = "Team A"
$ws.Range("G7").Value = "Team A"

# Leave the sheet's active selection on G7, matching the last cell touched.
$ws.Activate()
$ws.Range("G7").Select()
